$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. '24.455.94', '313.82').
# Prefix with an apostrophe so Excel stores them as text, matching the
# original inline-string cell type instead of coercing to a number.

$ws.Range("D2").Value = "'24.455.94"
$ws.Range("E2").Value = "  -1.06%  "
$ws.Range("D3").Value = "'1.689.36"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'313.82"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.3876"
$ws.Range("E7").Value = "  -2.06%  "
$ws.Range("D8").Value = "'0.4016"
$ws.Range("E8").Value = "  -1.47%  "
$ws.Range("D9").Value = "'1.489"
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("D10").Value = "'1.005"
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("D11").Value = "'52.41"
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("D12").Value = "'0.08739"
$ws.Range("E12").Value = "  -2.20%  "
$ws.Range("D13").Value = "'25.00"
$ws.Range("E13").Value = "  +6.24%  "
$ws.Range("D14").Value = "'7.490"
$ws.Range("E14").Value = "  +2.95%  "
$ws.Range("D15").Value = "'8.000"
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("D16").Value = "'0.00001344"
$ws.Range("E16").Value = "  +1.37%  "
$ws.Range("D17").Value = "'1.684.66"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").Value = "'98.16"
$ws.Range("E18").Value = "  -1.97%  "
$ws.Range("D19").Value = "'0.07076"
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").Value = "'19.96"
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("D21").Value = "'7.240"
$ws.Range("E21").Value = "  +3.46%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "'14.24"
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("D24").Value = "'24.447.86"
$ws.Range("E24").Value = "  -0.99%  "
$ws.Range("D25").Value = "'2.350"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("D26").Value = "'2.943"
$ws.Range("E26").Value = "  -10.07%  "
$ws.Range("D27").Value = "'22.66"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").Value = "'162.72"
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("D29").Value = "'8.603"
$ws.Range("E29").Value = "  +14.28%  "
$ws.Range("D30").Value = "'136.67"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").Value = "'5.205"
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("D32").Value = "'1.870.18"
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("D33").Value = "'0.08793"
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("D34").Value = "'7.394"
$ws.Range("E34").Value = "  +4.54%  "
$ws.Range("D35").Value = "'1.029"
$ws.Range("E35").Value = "  -2.46%  "
$ws.Range("D36").Value = "'0.2826"
$ws.Range("E36").Value = "  +2.95%  "
$ws.Range("D37").Value = "'1.954"
$ws.Range("E37").Value = "  +4.04%  "
$ws.Range("D38").Value = "'0.02892"
$ws.Range("E38").Value = "  +5.91%  "
$ws.Range("D39").Value = "'10.73"
$ws.Range("E39").Value = "  -5.90%  "
$ws.Range("D40").Value = "'14.16"
$ws.Range("E40").Value = "  -2.25%  "
$ws.Range("D41").Value = "'0.09101"
$ws.Range("E41").Value = "  -1.49%  "
$ws.Range("D42").Value = "'0.7918"
$ws.Range("E42").Value = "  +3.28%  "
$ws.Range("D43").Value = "'1.453"
$ws.Range("E43").Value = "  -1.60%  "
$ws.Range("D44").Value = "'16.78"
$ws.Range("E44").Value = "  +3.83%  "
$ws.Range("D45").Value = "'0.7230"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").Value = "'2.592"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").Value = "'4.198"
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("E49").Value = "  +1.43%  "
$ws.Range("D50").Value = "'138.12"
$ws.Range("E50").Value = "  -1.63%  "
$ws.Range("D51").Value = "'0.08018"
$ws.Range("E51").Value = "  +0.48%  "
